$wb = $excel.ActiveWorkbook

# Duplicate the "Portugal" sheet to create the new "Italy" sheet at the end
$portugal = $wb.Worksheets.Item("Portugal")
$portugal.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"

# Rows 3-4 on the new sheet go back to the default row height (no explicit ht);
# row 5 keeps its taller (28.8) height.
$italy.Rows("3:4").AutoFit()

# B4 gets new content with default (no) styling
$italy.Range("B4").Style = "Normal"
$italy.Range("B4").Value = "NGC-3145/T2159"

# Update the selection on the new Italy sheet (it is already the active tab)
$italy.Range("B2").Select()

# Swiss is no longer the active tab; its selection moved to B28
$swiss = $wb.Worksheets.Item("Swiss")
$swiss.Activate()
$swiss.Range("B28").Select()

# Portugal's selection becomes a full-sheet selection
$portugal.Activate()
$portugal.Cells.Select()

# Restore Italy as the active/selected tab
$italy.Activate()
